$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must remain plain text even when it looks numeric
# (matches the source data, which stores every Price cell as inline text).
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "90.809.82"
$ws.Range("E2").Value = "  +0.85%  "
$ws.Range("D3").Value = "3.177.67"
$ws.Range("E3").Value = "  +2.45%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.15%  "
Set-TextValue "D5" "218.75"
$ws.Range("E5").Value = "  +2.05%  "
Set-TextValue "D6" "626.51"
$ws.Range("E6").Value = "  +1.37%  "
Set-TextValue "D7" "1.12"
$ws.Range("E7").Value = "  +26.04%  "
Set-TextValue "D8" "0.371"
$ws.Range("E8").Value = "  +0.25%  "
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "3.173.09"
$ws.Range("E10").Value = "  +2.31%  "
Set-TextValue "D11" "0.757"
$ws.Range("E11").Value = "  +13.88%  "
Set-TextValue "D12" "0.200"
$ws.Range("E12").Value = "  +6.34%  "
Set-TextValue "D13" "0.0000249"
$ws.Range("E13").Value = "  +2.91%  "
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D14" "5.66"
$ws.Range("E14").Value = "  +5.10%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextValue "D15" "35.24"
$ws.Range("E15").Value = "  +7.87%  "
$ws.Range("D16").Value = "90.663.17"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "3.748.43"
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "3.125.97"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("E19").Value = "  +10.26%  "
$ws.Range("E20").Value = "  +3.27%  "
Set-TextValue "D21" "14.38"
$ws.Range("E21").Value = "  +5.55%  "
Set-TextValue "D22" "446.81"
$ws.Range("E22").Value = "  +2.82%  "
Set-TextValue "D23" "8.97"
$ws.Range("E23").Value = "  +8.46%  "
Set-TextValue "D24" "5.23"
$ws.Range("E24").Value = "  +3.91%  "
Set-TextValue "D25" "6.01"
$ws.Range("E25").Value = "  +8.71%  "
Set-TextValue "D26" "87.87"
$ws.Range("E26").Value = "  +1.98%  "
Set-TextValue "D27" "12.33"
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("D28").Value = "3.341.63"
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("B30").Value = "Cronos"
$ws.Range("C30").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D30" "0.162"
$ws.Range("E30").Value = "  -2.24%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D31" "9.28"
$ws.Range("E31").Value = "  +13.06%  "
Set-TextValue "D32" "1.01"
$ws.Range("E32").Value = "  -7.59%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D33" "26.03"
$ws.Range("E33").Value = "  +13.12%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D34" "527.26"
$ws.Range("E34").Value = "  +2.78%  "
Set-TextValue "D35" "3.74"
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D36" "0.145"
$ws.Range("E36").Value = "  +9.03%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
Set-TextValue "D37" "7.03"
$ws.Range("E37").Value = "  +3.88%  "
$ws.Range("E38").Value = "  +5.94%  "
$ws.Range("E39").Value = "  +4.41%  "
Set-TextValue "D40" "0.174"
$ws.Range("E40").Value = "  +21.20%  "
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D42" "0.0855"
$ws.Range("E42").Value = "  +23.89%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D43" "1.00"
$ws.Range("E43").Value = "  -0.23%  "
Set-TextValue "D44" "0.416"
$ws.Range("E44").Value = "  +11.35%  "
Set-TextValue "D45" "1.95"
$ws.Range("E45").Value = "  +5.22%  "
$ws.Range("E46").Value = "  +0.00%  "
Set-TextValue "D47" "148.41"
$ws.Range("E47").Value = "  +0.98%  "
Set-TextValue "D48" "1.35"
$ws.Range("E48").Value = "  +9.69%  "
Set-TextValue "D49" "44.08"
$ws.Range("E49").Value = "  +1.14%  "
Set-TextValue "D50" "4.40"
$ws.Range("E50").Value = "  +6.92%  "
Set-TextValue "D51" "0.657"
$ws.Range("E51").Value = "  +10.71%  "
